$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to be treated as text so numeric-looking values
# (e.g. "581.68") are not coerced into floating point numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "69.300.81"
$ws.Range("E2").Value = "  -2.55%  "
$ws.Range("D3").Value = "3.534.10"
$ws.Range("E3").Value = "  -4.32%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "581.68"
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("D6").Value = "171.99"
$ws.Range("E6").Value = "  -3.53%  "
$ws.Range("D7").Value = "3.527.56"
$ws.Range("E7").Value = "  -4.25%  "
$ws.Range("D8").Value = "0.609"
$ws.Range("E8").Value = "  -1.26%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("E10").Value = "  -5.30%  "
$ws.Range("D11").Value = "6.71"
$ws.Range("E11").Value = "  -1.53%  "
$ws.Range("D12").Value = "0.585"
$ws.Range("E12").Value = "  -4.41%  "
$ws.Range("D13").Value = "47.44"
$ws.Range("E13").Value = "  -3.51%  "
$ws.Range("E14").Value = "  -4.71%  "
$ws.Range("D15").Value = "4.093.42"
$ws.Range("E15").Value = "  -4.55%  "
$ws.Range("D16").Value = "8.55"
$ws.Range("E16").Value = "  -5.70%  "
$ws.Range("D17").Value = "627.54"
$ws.Range("E17").Value = "  -7.52%  "
$ws.Range("D18").Value = "3.544.39"
$ws.Range("E18").Value = "  -3.99%  "
$ws.Range("D19").Value = "69.280.17"
$ws.Range("E19").Value = "  -2.76%  "
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").Value = "17.52"
$ws.Range("E21").Value = "  -2.70%  "
$ws.Range("D22").Value = "11.22"
$ws.Range("E22").Value = "  -3.17%  "
$ws.Range("D23").Value = "0.889"
$ws.Range("E23").Value = "  -5.78%  "
$ws.Range("D24").Value = "15.97"
$ws.Range("E24").Value = "  -8.33%  "
$ws.Range("D25").Value = "97.66"
$ws.Range("E25").Value = "  -4.30%  "
$ws.Range("E26").Value = "  -4.14%  "
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").Value = "2.64"
$ws.Range("E28").Value = "  -7.05%  "
$ws.Range("D29").Value = "9.31"
$ws.Range("E29").Value = "  -9.10%  "
$ws.Range("D30").Value = "32.89"
$ws.Range("E30").Value = "  -6.49%  "
$ws.Range("E31").Value = "  -7.64%  "
$ws.Range("D32").Value = "8.56"
$ws.Range("E32").Value = "  -6.59%  "
$ws.Range("E33").Value = "  -6.66%  "
$ws.Range("D34").Value = "7.00"
$ws.Range("E34").Value = "  -7.67%  "
$ws.Range("D35").Value = "633.36"
$ws.Range("E35").Value = "  +8.56%  "
$ws.Range("D36").Value = "10.78"
$ws.Range("E36").Value = "  -3.94%  "
$ws.Range("D37").Value = "3.50"
$ws.Range("E37").Value = "  -14.25%  "
$ws.Range("E38").Value = "  -5.33%  "
$ws.Range("D39").Value = "57.45"
$ws.Range("E39").Value = "  -2.46%  "
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("D41").Value = "0.0451"
$ws.Range("E41").Value = "  -2.59%  "
$ws.Range("E42").Value = "  -5.17%  "
$ws.Range("D43").Value = "3.390.92"
$ws.Range("E43").Value = "  -6.59%  "
$ws.Range("D44").Value = "0.329"
$ws.Range("E44").Value = "  -6.46%  "
$ws.Range("D45").Value = "32.93"
$ws.Range("E45").Value = "  -6.86%  "
$ws.Range("D46").Value = "0.0₃0701"
$ws.Range("E46").Value = "  -8.75%  "
$ws.Range("D47").Value = "2.56"
$ws.Range("E47").Value = "  -7.50%  "
$ws.Range("E48").Value = "  -4.55%  "
$ws.Range("E49").Value = "  -2.71%  "
$ws.Range("D50").Value = "5.70"
$ws.Range("E50").Value = "  +14.09%  "
$ws.Range("D51").Value = "132.11"
$ws.Range("E51").Value = "  -2.07%  "

# Restore the default (unstyled) cell style so formatting matches the
# original workbook, while keeping the values stored as text.
$ws.Range("D2:E51").Style = "Normal"
